# Refresh the "抓取标题" (scraped-titles) sheet with a new batch of
# scraped news rows (rows 2-21). Header row (row 1: Title/url/content/
# img_url/img_content) and the other sheets are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("抓取标题")

$ws.Cells.Item(2, 1).Value = '是捍卫言论自由还是整顿网络暴力'
$ws.Cells.Item(2, 2).Value = 'https://www.jinse.com/bitcoin/276430.html'
$ws.Cells.Item(2, 3).Value = '是捍卫言论自由还是整顿网络暴力'
$ws.Cells.Item(2, 4).Value = 'https://img.jinse.com/1345469_image1.png'

$ws.Cells.Item(3, 1).Value = '区块链十年 | 第1.2章  2008 比特币白皮书'
$ws.Cells.Item(3, 2).Value = 'https://www.jinse.com/member/208800'
$ws.Cells.Item(3, 3).Value = '在世界的某个角落里，中本聪开发的第一版比特币客户端也即将完成。'
$ws.Cells.Item(3, 4).Value = 'https://img.jinse.com/1345449_image1.png'

$ws.Cells.Item(4, 1).Value = '懂个币 BSV扛不住了？算力战已让矿工损失数百万美元收入'
$ws.Cells.Item(4, 2).Value = 'https://www.jinse.com/bitcoin/276429.html'
$ws.Cells.Item(4, 3).Value = '算力战让矿工损失数百万美元收入。'
$ws.Cells.Item(4, 4).Value = 'https://img.jinse.com/1343898_image1.png'

$ws.Cells.Item(5, 1).Value = '4分钟课堂 | 经济泡沫是如何产生的？'
$ws.Cells.Item(5, 2).Value = 'https://www.jinse.com/member/190014'
$ws.Cells.Item(5, 3).Value = '自17世纪荷兰“郁金香狂热”事件成为第一场经济泡沫后，市场涌现出房地产泡沫 、股市泡沫、互联网泡沫、甚至比特币泡沫等各种经济泡沫。经济学家将经济泡沫分为五个阶段'
$ws.Cells.Item(5, 4).Value = 'https://img.jinse.com/1345452_image1.png'

$ws.Cells.Item(6, 1).Value = '区块链十年 | 第1.1章  2008 世界金融危机'
$ws.Cells.Item(6, 2).Value = 'https://www.jinse.com/bitcoin/276173.html'
$ws.Cells.Item(6, 3).Value = '面对金融危机，一位神秘人思考着自己的解决方案......'
$ws.Cells.Item(6, 4).Value = 'https://img.jinse.com/1345449_image1.png'

$ws.Cells.Item(7, 1).Value = '2018.11.21 今日早盘短线策略 遵守规则敬畏市场'
$ws.Cells.Item(7, 2).Value = 'https://www.jinse.com/member/211555'
$ws.Cells.Item(7, 3).Value = '2018.11.21 早间策略一览'
$ws.Cells.Item(7, 4).Value = 'https://img.jinse.com/1298737_image1.png'

$ws.Cells.Item(8, 1).Value = '哈希派：多头反弹软绵绵、但空头也想休息一下了'
$ws.Cells.Item(8, 2).Value = 'https://www.jinse.com/bitcoin/276427.html'
$ws.Cells.Item(8, 3).Value = '今日行情分析币种：BTC / ETH / EOS / BCH / XRP / LTC'
$ws.Cells.Item(8, 4).Value = 'https://img.jinse.com/1330734_image1.png'

$ws.Cells.Item(9, 1).Value = '万点财经熊市专栏|陶太胜的熊市生存法则'
$ws.Cells.Item(9, 2).Value = 'https://www.jinse.com/member/202987'
$ws.Cells.Item(9, 3).Value = '我本尘埃，因花而来'
$ws.Cells.Item(9, 4).Value = 'https://img.jinse.com/1345385_image1.png'

$ws.Cells.Item(10, 1).Value = '加密货币市值24小时蒸发2600亿“避险属性”丧失共识肢解？'
$ws.Cells.Item(10, 2).Value = 'https://www.jinse.com/bitcoin/276426.html'
$ws.Cells.Item(10, 3).Value = '这场暴跌背后，有更深层的原因和更深远的影响。'
$ws.Cells.Item(10, 4).Value = 'https://img.jinse.com/1343868_image1.png'

$ws.Cells.Item(11, 1).Value = '走进区块链企业 | 赵伟：冬天很冷但我们很感谢这个冬天'
$ws.Cells.Item(11, 2).Value = 'https://www.jinse.com/member/190014'
$ws.Cells.Item(11, 3).Value = '这是一家走进去，便想细细去探索一番的企业'
$ws.Cells.Item(11, 4).Value = 'https://img.jinse.com/1343833_image1.png'

$ws.Cells.Item(12, 1).Value = '人工智能与区块链的边界逐渐模糊，软硬件将如何结合？'
$ws.Cells.Item(12, 2).Value = 'https://www.jinse.com/bitcoin/276425.html'
$ws.Cells.Item(12, 3).Value = '2018年11月18日，2018比原链全球开发者大会在杭州国际博览中心（G20会馆）继续进行，这是杭州第一次由开源组织举办的技术型峰会，也是杭州被誉为区块链之城以来规模最大的一场区块链开发者'
$ws.Cells.Item(12, 4).Value = 'https://img.jinse.com/1345437_image1.png'

$ws.Cells.Item(13, 1).Value = 'OGC第一款DAPP即将震撼来袭'
$ws.Cells.Item(13, 2).Value = 'https://www.jinse.com/member/208546'
$ws.Cells.Item(13, 3).Value = 'DAPP是OGC技术应用落地的第一步，也是其对全球网游行业颠覆之路的正式实施。'
$ws.Cells.Item(13, 4).Value = 'https://img.jinse.com/1345335_image1.png'

$ws.Cells.Item(14, 1).Value = '兄弟们秋裤不顶事了棉裤穿上吧'
$ws.Cells.Item(14, 2).Value = 'https://www.jinse.com/bitcoin/276424.html'
$ws.Cells.Item(14, 3).Value = '矿机论斤卖。'
$ws.Cells.Item(14, 4).Value = 'https://img.jinse.com/1343806_image1.png'

$ws.Cells.Item(15, 1).Value = 'AiBiChain创始人刘靖中入围“2018中国创新榜样人物”'
$ws.Cells.Item(15, 2).Value = 'https://www.jinse.com/member/26030'
$ws.Cells.Item(15, 3).Value = '作为数字技术的前沿，近年来区块链技术在许多国家受到越来越多的关注。在中国，区块链与其他几个重要新技术一并被列为《“十三五”国家信息化规划》中的“重大任务和重点'
$ws.Cells.Item(15, 4).Value = 'https://img.jinse.com/1345305_image1.png'

$ws.Cells.Item(16, 1).Value = '网信办主管媒体针对区块链连发8文：肯定Token价值 落地应用是重点'
$ws.Cells.Item(16, 2).Value = 'https://www.jinse.com/bitcoin/276423.html'
$ws.Cells.Item(16, 3).Value = '结合区块链目前“遇冷”的行业背景，或许可从这几篇组稿以一窥网信办等国家机构现阶段对于区块链的立场和态度。'
$ws.Cells.Item(16, 4).Value = 'https://img.jinse.com/1345340_image1.png'

$ws.Cells.Item(17, 1).Value = '“大崩盘”过后策略抄底的机会来了'
$ws.Cells.Item(17, 2).Value = 'https://www.jinse.com/member/142128'
$ws.Cells.Item(17, 3).Value = '捡垃圾的时候到了'
$ws.Cells.Item(17, 4).Value = 'https://img.jinse.com/1343795_image1.png'

$ws.Cells.Item(18, 1).Value = '谁也没有想到 区块链还能用来熬鸡汤'
$ws.Cells.Item(18, 2).Value = 'https://www.jinse.com/bitcoin/276169.html'
$ws.Cells.Item(18, 3).Value = '希望这是币圈人熬的最后一碗“鸡汤”。'
$ws.Cells.Item(18, 4).Value = 'https://img.jinse.com/1345343_image1.png'

$ws.Cells.Item(19, 1).Value = '币易Coinyee.io交易所1周年 百万Token大派送 100%有礼'
$ws.Cells.Item(19, 2).Value = 'https://www.jinse.com/member/126549'
$ws.Cells.Item(19, 3).Value = '为感谢广大用户朋友们一直以来的陪伴与支持，币易Coinyee.io特借助1周年生日的机会，举办了一场价值百万元Token好礼派送活动，100%有礼！'
$ws.Cells.Item(19, 4).Value = 'https://img.jinse.com/1345140_image1.png'

$ws.Cells.Item(20, 1).Value = '英雄会聚 恰逢其时 元界CTO陈浩受邀参加Money20/20峰会圆桌会议 共话区块链最新趋势'
$ws.Cells.Item(20, 2).Value = 'https://www.jinse.com/bitcoin/276163.html'
$ws.Cells.Item(20, 3).Value = '元界CTO陈浩受Money20/20主办方邀请参与圆桌会谈环节与全球“重量级”演讲嘉宾就如何打造金融区块链进行了深入的交流探讨。'
$ws.Cells.Item(20, 4).Value = 'https://img.jinse.com/1345331_image1.png'

$ws.Cells.Item(21, 1).Value = 'NVIDIA Q3盈利报告印证了加密货币采矿业务正在枯竭'
$ws.Cells.Item(21, 2).Value = 'https://www.jinse.com/member/214042'
$ws.Cells.Item(21, 3).Value = '美国芯片龙头科技公司NVIDIA在2018年11月15日发布了2018年第三季度财报。该报告指出'
$ws.Cells.Item(21, 4).Value = 'https://img.jinse.com/1343723_image1.png'
